$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "1.000") are not
# auto-converted to numbers by Excel when assigned via .Value
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.030.46'
$ws.Range("E2").Value = '  -1.48%  '

$ws.Range("D3").Value = '1.827.16'
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.44%  '

$ws.Range("D5").Value = '311.51'
$ws.Range("E5").Value = '  -1.17%  '

$ws.Range("E6").Value = '  -0.28%  '

$ws.Range("D7").Value = '0.4257'
$ws.Range("E7").Value = '  -0.94%  '

$ws.Range("D8").Value = '0.3660'
$ws.Range("E8").Value = '  -1.15%  '

$ws.Range("D9").Value = '0.07248'
$ws.Range("E9").Value = '  -0.17%  '

$ws.Range("D10").Value = '0.8440'
$ws.Range("E10").Value = '  -2.62%  '

$ws.Range("D11").Value = '20.64'
$ws.Range("E11").Value = '  -2.51%  '

$ws.Range("D12").Value = '1.821.19'
$ws.Range("E12").Value = '  -1.07%  '

$ws.Range("D13").Value = '6.661'
$ws.Range("E13").Value = '  -0.54%  '

$ws.Range("D14").Value = '0.07051'
$ws.Range("E14").Value = '  -0.23%  '

$ws.Range("D15").Value = '5.291'
$ws.Range("E15").Value = '  -1.28%  '

$ws.Range("D16").Value = '89.69'
$ws.Range("E16").Value = '  +1.97%  '

$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.30%  '

$ws.Range("D18").Value = '0.000008770'
$ws.Range("E18").Value = '  -1.62%  '

$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("E20").Value = '  -2.28%  '

$ws.Range("D21").Value = '27.131.40'
$ws.Range("E21").Value = '  -1.11%  '

$ws.Range("D22").Value = '5.138'
$ws.Range("E22").Value = '  -0.55%  '

$ws.Range("D23").Value = '10.85'
$ws.Range("E23").Value = '  -0.50%  '

$ws.Range("D24").Value = '2.054.95'
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("D25").Value = '1.983'
$ws.Range("E25").Value = '  -0.93%  '

$ws.Range("D26").Value = '151.20'
$ws.Range("E26").Value = '  -1.48%  '

$ws.Range("D27").Value = '2.221'

$ws.Range("D28").Value = '18.27'
$ws.Range("E28").Value = '  -0.97%  '

$ws.Range("D29").Value = '5.231'
$ws.Range("E29").Value = '  -1.10%  '

$ws.Range("D30").Value = '116.98'
$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("D31").Value = '0.08708'
$ws.Range("E31").Value = '  -1.70%  '

$ws.Range("D32").Value = '1.179'
$ws.Range("E32").Value = '  -2.59%  '

$ws.Range("D33").Value = '0.7399'
$ws.Range("E33").Value = '  -3.45%  '

$ws.Range("D34").Value = '2.901'
$ws.Range("E34").Value = '  -0.25%  '

$ws.Range("D35").Value = '4.422'
$ws.Range("E35").Value = '  -1.41%  '

$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("D37").Value = '1.093'
$ws.Range("E37").Value = '  -2.55%  '

$ws.Range("D38").Value = '0.01944'
$ws.Range("E38").Value = '  -0.84%  '

$ws.Range("D39").Value = '0.05224'
$ws.Range("E39").Value = '  -1.26%  '

$ws.Range("D40").Value = '7.234'
$ws.Range("E40").Value = '  +0.72%  '

$ws.Range("D41").Value = '2.870'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1691'
$ws.Range("E42").Value = '  +1.05%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.5123'
$ws.Range("E43").Value = '  +0.73%  '

$ws.Range("D44").Value = '8.569'
$ws.Range("E44").Value = '  -0.91%  '

$ws.Range("D45").Value = '10.57'
$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").Value = '0.4740'
$ws.Range("E46").Value = '  -0.16%  '

$ws.Range("D47").Value = '1.943'
$ws.Range("E47").Value = '  +6.14%  '

$ws.Range("D48").Value = '105.68'
$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("D49").Value = '0.9995'
$ws.Range("E49").Value = '  -0.32%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.659'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.06327'
$ws.Range("E51").Value = '  -1.53%  '

# Restore column D cell style back to the workbook default ("Normal") now that the
# text values are set, so no stray number-format style lingers on the cells.
$ws.Range("D2:D51").Style = "Normal"
